# Apply scheduled market-price / profit recalculation updates to each Leve-profit sheet.
# Values correspond to updated currentAveragePrice / LevePrice / LeveProfit columns (H:N)
# pulled by the scheduled runner for each class sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 20.944445
$ws.Range("I11").Value = 20.944445
$ws.Range("K11").Value = 20.944445
$ws.Range("M11").Value = 119.055555
$ws.Range("H52").Value = 440.625
$ws.Range("I52").Value = 440.83334
$ws.Range("J52").Value = 440
$ws.Range("K52").Value = 1322.50002
$ws.Range("L52").Value = 1320
$ws.Range("M52").Value = -1162.50002
$ws.Range("N52").Value = -1640
$ws.Range("H62").Value = 14650.111
$ws.Range("I62").Value = 10550.143
$ws.Range("J62").Value = 29000
$ws.Range("K62").Value = 10550.143
$ws.Range("L62").Value = 29000
$ws.Range("M62").Value = -9926.143
$ws.Range("N62").Value = -30248
$ws.Range("H65").Value = 14650.111
$ws.Range("I65").Value = 10550.143
$ws.Range("J65").Value = 29000
$ws.Range("K65").Value = 52750.715
$ws.Range("L65").Value = 145000
$ws.Range("M65").Value = -49630.715
$ws.Range("N65").Value = -151240
$ws.Range("H107").Value = 287.22223
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H112").Value = 1668.4222
$ws.Range("J112").Value = 1748.683
$ws.Range("L112").Value = 5246.049
$ws.Range("N112").Value = -7462.049
$ws.Range("H125").Value = 3038.5
$ws.Range("I125").Value = 3379.375
$ws.Range("J125").Value = 2765.8
$ws.Range("K125").Value = 30414.375
$ws.Range("L125").Value = 24892.2
$ws.Range("M125").Value = -27954.375
$ws.Range("N125").Value = -29812.2
$ws.Range("H129").Value = 1857.6522
$ws.Range("I129").Value = 1364.579
$ws.Range("J129").Value = 4199.75
$ws.Range("K129").Value = 4093.737
$ws.Range("L129").Value = 12599.25
$ws.Range("M129").Value = 906.2629999999999
$ws.Range("N129").Value = -22599.25
$ws.Range("H138").Value = 3700.5
$ws.Range("J138").Value = 4926.0435
$ws.Range("L138").Value = 14778.1305
$ws.Range("N138").Value = -25058.1305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2222.0322
$ws.Range("I2").Value = 2029.8889
$ws.Range("K2").Value = 2029.8889
$ws.Range("M2").Value = -1916.8889
$ws.Range("H31").Value = 18316.334
$ws.Range("I31").Value = 18316.334
$ws.Range("K31").Value = 18316.334
$ws.Range("M31").Value = -18022.334
$ws.Range("H32").Value = 7578932
$ws.Range("I32").Value = 7815742.5
$ws.Range("J32").Value = 1004.5
$ws.Range("K32").Value = 7815742.5
$ws.Range("L32").Value = 1004.5
$ws.Range("M32").Value = -7815455.5
$ws.Range("N32").Value = -1578.5
$ws.Range("H37").Value = 47549.4
$ws.Range("I37").Value = 11199.6
$ws.Range("J37").Value = 59666
$ws.Range("K37").Value = 11199.6
$ws.Range("L37").Value = 59666
$ws.Range("M37").Value = -10926.6
$ws.Range("N37").Value = -60212
$ws.Range("H45").Value = 4067.6155
$ws.Range("I45").Value = 3973.1875
$ws.Range("J45").Value = 4218.7
$ws.Range("K45").Value = 3973.1875
$ws.Range("L45").Value = 4218.7
$ws.Range("M45").Value = -3596.1875
$ws.Range("N45").Value = -4972.7
$ws.Range("H61").Value = 2781294.8
$ws.Range("J61").Value = 2105.5
$ws.Range("L61").Value = 2105.5
$ws.Range("N61").Value = -2529.5
$ws.Range("H110").Value = 557.6667
$ws.Range("I110").Value = 557.6667
$ws.Range("K110").Value = 557.6667
$ws.Range("M110").Value = 1487.3333
$ws.Range("H116").Value = 2222.0322
$ws.Range("I116").Value = 2029.8889
$ws.Range("K116").Value = 2029.8889
$ws.Range("M116").Value = 264.1111000000001
$ws.Range("H132").Value = 805443.3
$ws.Range("I132").Value = 1013136.8
$ws.Range("K132").Value = 3039410.4
$ws.Range("M132").Value = -3036880.4
$ws.Range("H136").Value = 2781294.8
$ws.Range("J136").Value = 2105.5
$ws.Range("L136").Value = 6316.5
$ws.Range("N136").Value = -11416.5
$ws.Range("H138").Value = 106249.75
$ws.Range("J138").Value = 106249.75
$ws.Range("L138").Value = 106249.75
$ws.Range("N138").Value = -116529.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 78949.5
$ws.Range("J2").Value = 78949.5
$ws.Range("L2").Value = 78949.5
$ws.Range("N2").Value = -79175.5
$ws.Range("H3").Value = 2222.0322
$ws.Range("I3").Value = 2029.8889
$ws.Range("K3").Value = 2029.8889
$ws.Range("M3").Value = -1915.8889
$ws.Range("H134").Value = 462461.97
$ws.Range("I134").Value = 540783.75
$ws.Range("K134").Value = 1622351.25
$ws.Range("M134").Value = -1619816.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 74154.8
$ws.Range("J70").Value = 74154.8
$ws.Range("L70").Value = 74154.8
$ws.Range("N70").Value = -74784.8
$ws.Range("H73").Value = 74154.8
$ws.Range("J73").Value = 74154.8
$ws.Range("L73").Value = 74154.8
$ws.Range("N73").Value = -76338.8
$ws.Range("H107").Value = 1065.3
$ws.Range("I107").Value = 1175.2858
$ws.Range("J107").Value = 808.6667
$ws.Range("K107").Value = 1175.2858
$ws.Range("L107").Value = 808.6667
$ws.Range("M107").Value = 744.7141999999999
$ws.Range("N107").Value = -4648.6667
$ws.Range("H132").Value = 5216806.5
$ws.Range("I132").Value = 8798.325999999999
$ws.Range("K132").Value = 26394.978
$ws.Range("M132").Value = -23864.978
$ws.Range("H134").Value = 2305.7144
$ws.Range("I134").Value = 2328.1538
$ws.Range("K134").Value = 6984.4614
$ws.Range("M134").Value = -4449.4614
$ws.Range("H140").Value = 73331.336
$ws.Range("J140").Value = 73331.336
$ws.Range("L140").Value = 73331.336
$ws.Range("N140").Value = -83691.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 14921.2
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 14921.2
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 44763.60000000001
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -46261.60000000001
$ws.Range("H66").Value = 14921.2
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 14921.2
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 134290.8
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -141778.8
$ws.Range("H74").Value = 7197.8335
$ws.Range("I74").Value = 4499.6665
$ws.Range("J74").Value = 9896
$ws.Range("K74").Value = 13498.9995
$ws.Range("L74").Value = 29688
$ws.Range("M74").Value = -12437.9995
$ws.Range("N74").Value = -31810
$ws.Range("H77").Value = 7197.8335
$ws.Range("I77").Value = 4499.6665
$ws.Range("J77").Value = 9896
$ws.Range("K77").Value = 40496.9985
$ws.Range("L77").Value = 89064
$ws.Range("M77").Value = -35192.9985
$ws.Range("N77").Value = -99672
$ws.Range("H93").Value = 4170.091
$ws.Range("J93").Value = 5109
$ws.Range("L93").Value = 15327
$ws.Range("N93").Value = -19071
$ws.Range("H110").Value = 17495.5
$ws.Range("I110").Value = 2426.5
$ws.Range("J110").Value = 25030
$ws.Range("K110").Value = 7279.5
$ws.Range("L110").Value = 75090
$ws.Range("M110").Value = -3189.5
$ws.Range("N110").Value = -83270
$ws.Range("H119").Value = 2999
$ws.Range("I119").Value = 1998.5
$ws.Range("K119").Value = 5995.5
$ws.Range("M119").Value = -1157.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 3500
$ws.Range("I4").Value = 3500
$ws.Range("K4").Value = 3500
$ws.Range("M4").Value = -3388
$ws.Range("H11").Value = 4576.75
$ws.Range("I11").Value = 5102.3335
$ws.Range("K11").Value = 5102.3335
$ws.Range("M11").Value = -4963.3335
$ws.Range("H17").Value = 8249.333000000001
$ws.Range("I17").Value = 5665.3335
$ws.Range("J17").Value = 10833.333
$ws.Range("K17").Value = 5665.3335
$ws.Range("L17").Value = 10833.333
$ws.Range("M17").Value = -5497.3335
$ws.Range("N17").Value = -11169.333
$ws.Range("H113").Value = 3260.0715
$ws.Range("I113").Value = 1488.5
$ws.Range("J113").Value = 3968.7
$ws.Range("K113").Value = 1488.5
$ws.Range("L113").Value = 3968.7
$ws.Range("M113").Value = 681.5
$ws.Range("N113").Value = -8308.700000000001
$ws.Range("H122").Value = 4320.129
$ws.Range("I122").Value = 4416.4
$ws.Range("J122").Value = 3919
$ws.Range("K122").Value = 13249.2
$ws.Range("L122").Value = 11757
$ws.Range("M122").Value = -10799.2
$ws.Range("N122").Value = -16657
$ws.Range("H132").Value = 929690.1
$ws.Range("I132").Value = 1507898.9
$ws.Range("J132").Value = 4556.2
$ws.Range("K132").Value = 4523696.699999999
$ws.Range("L132").Value = 13668.6
$ws.Range("M132").Value = -4521166.699999999
$ws.Range("N132").Value = -18728.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2923.9048
$ws.Range("I61").Value = 1632.2307
$ws.Range("J61").Value = 5022.875
$ws.Range("K61").Value = 1632.2307
$ws.Range("L61").Value = 5022.875
$ws.Range("M61").Value = -1430.2307
$ws.Range("N61").Value = -5426.875
$ws.Range("H113").Value = 2923.9048
$ws.Range("I113").Value = 1632.2307
$ws.Range("J113").Value = 5022.875
$ws.Range("K113").Value = 1632.2307
$ws.Range("L113").Value = 5022.875
$ws.Range("M113").Value = 537.7692999999999
$ws.Range("N113").Value = -9362.875
$ws.Range("H122").Value = 3592.6155
$ws.Range("I122").Value = 3352.6206
$ws.Range("J122").Value = 4288.6
$ws.Range("K122").Value = 10057.8618
$ws.Range("L122").Value = 12865.8
$ws.Range("M122").Value = -7607.861800000001
$ws.Range("N122").Value = -17765.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 160023.88
$ws.Range("J74").Value = 184865.5
$ws.Range("L74").Value = 184865.5
$ws.Range("N74").Value = -186737.5
$ws.Range("H77").Value = 160023.88
$ws.Range("J77").Value = 184865.5
$ws.Range("L77").Value = 554596.5
$ws.Range("N77").Value = -563956.5
$ws.Range("H107").Value = 2525.4614
$ws.Range("I107").Value = 1081.0834
$ws.Range("J107").Value = 3763.5
$ws.Range("K107").Value = 3243.2502
$ws.Range("L107").Value = 11290.5
$ws.Range("M107").Value = -1323.2502
$ws.Range("N107").Value = -15130.5
$ws.Range("H113").Value = 2639.3635
$ws.Range("I113").Value = 1267.0555
$ws.Range("J113").Value = 4286.1333
$ws.Range("K113").Value = 3801.1665
$ws.Range("L113").Value = 12858.3999
$ws.Range("M113").Value = -1631.1665
$ws.Range("N113").Value = -17198.3999
